$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @(
    "-4,-4",
    "9,0",
    "-10,-2",
    "-3,5",
    "-6,-8",
    "3,-9",
    "8,4",
    "-8,2",
    "2,8",
    "-9,-1",
    "2,5",
    "-4,-1",
    "-7,-9",
    "9,-4",
    "7,1",
    "8,-6",
    "-5,4",
    "8,2",
    "5,0",
    "-4,-6"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}
